# Updated descriptions of sponsorship levels
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Gold Sponsor (row 6) benefits description ---
$ws.Range("E6").Value2 = "1. Mention of your organization in emails to event attendees`n2. Branded swag materials you provide included in attendee kit`n3. Verbal acknowledgment of your sponsorship at the event during Lunch`n4.  Display of your organization logo in our website Sponsors page`n5. Mention in at least 5 social media posts`n6. Preferred location for your booth at the conference`n8. Display of your logo in Eventbrite ticketing page`n"

# --- Silver Sponsor (row 7) benefits description ---
$ws.Range("E7").Value2 = "1. Mention of your organization in emails to event attendees `n2. Verbal acknowledgment of your sponsorship at the event during Lunch`n3. Display of your organization logo in our website Sponsors page`n4. Mention in at least 2 social media posts`n5. Space for your organization's booth at the conference"

# --- Bronze Sponsor (row 8) benefits description ---
$ws.Range("E8").Value2 = "1. Display of your organization logo in our  website Sponsors page`n2. Mention in at least 1 social media post"

# --- Speaker Travel Sponsor (row 9) benefits description ---
$ws.Range("E9").Value2 = "1.  Display of your organization logo in our website Additional Sponsors Section"

# --- Lunch Sponsor (row 10) benefits description ---
$ws.Range("E10").Value2 = "1. Verbal acknowledgment of your sponsorship at the event during lunch`n2. Display of your organization logo in our website Sponsors page`n3. Display of your organization logo in presentation deck at lunch`n"

# --- Speaker Dinner Sponsor (row 11) benefits description ---
$ws.Range("E11").Value2 = "1. Verbal acknowledgment of your sponsorship during Speaker Dinner`n2. Display of your organization logo in our website Sponsors section`n3. Display of your organization logo in presentation deck in Lunch Room"

# --- Row heights re-fitted to the new (shorter) text ---
$ws.Rows.Item(6).RowHeight = 120
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 45

# --- Restore the view/selection state left by the editor ---
$ws.Range("E18").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
